# The source data for this "Coliflor" sheet gained one more weekly record.
# A brand-new observation (Fecha 44858) is inserted as the new row 306,
# pushing the previously-existing rows 306..330 down to 307..331.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 306 (shifts rows 306:330 down to 307:331).
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row 306 with the new record.
$ws.Cells.Item(306, 1).Value  = 7
$ws.Cells.Item(306, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(306, 3).Value  = "Ñuble"
$ws.Cells.Item(306, 4).Value  = 44858
$ws.Cells.Item(306, 5).Value  = 16
$ws.Cells.Item(306, 6).Value  = 100112008
$ws.Cells.Item(306, 7).Value  = "Coliflor"
$ws.Cells.Item(306, 8).Value  = "Sin especificar"
$ws.Cells.Item(306, 9).Value  = "Primera"
$ws.Cells.Item(306, 10).Value = 400
$ws.Cells.Item(306, 11).Value = 1000
$ws.Cells.Item(306, 12).Value = 1200
$ws.Cells.Item(306, 13).Value = 1100
$ws.Cells.Item(306, 14).Value = "$/unidad"
$ws.Cells.Item(306, 15).Value = "Región del Maule"
$ws.Cells.Item(306, 16).Value = 1100
$ws.Cells.Item(306, 17).Value = 1
$ws.Cells.Item(306, 18).Value = "Hortaliza"
